$d = $word.ActiveDocument

# Use Find.Execute (not a bare Range.Text assignment) for every text swap so
# that xml:space="preserve" on the run is recomputed from the new text
# rather than being left over from whatever the original run had.
#
# Note: Word's Find/Replace treats "^" in the ReplaceWith text as the start
# of a special code (^t = tab, ^p = paragraph mark, ...), so a literal caret
# must be escaped as "^^" in the replacement text (e.g. "r^^T" -> "r^T").

# 1) Date in the header line: 05.11.24 -> 04.11.24
$d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק -05.11.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק -04.11.24: ⚡️🚀", 2)

# 2) Title
$d.Content.Find.Execute("RETHINKING SOFTMAX: SELF-ATTENTION WITH POLYNOMIAL ACTIVATIONS", $true, $false, $false, $false, $false, $true, 1, $false, "Refusal in Language Models Is Mediated by a Single Direction", 2)

# 3) First body paragraph
$d.Content.Find.Execute("מאמר די לא רגיל והוא מדבר על חלופה פוטנציאלית של מנגנון ה-attention שאנו כה אוהבים בטרנספורמים. אתם בטח זוכרים שמשקלי attention בשנאים מחושבים עם softmax שהוא מנרמל וקטורי משקלים לנורמה 1 ובנוסף כל רכיביו הינם בין 0 ל- 1 כלומר הוא מהווה התפלגות הסתברותית. המחברים טוענים שתכונות אלו של המשקלים לא קריטיות לפונקציונאליות של השנאים ומציעים להחליף אותם בקרנל אחר שהוא פולינומיאלי כפי שאתם בטח ניחשתם מהשם של המאמר.", $true, $false, $false, $false, $false, $true, 1, $false, "מאמר מעניין החוקר איך ניתן לגרום למודל שפה לתת תשובות רצויות יותר ורצויות פחות. מתברר שאפשר לגרום למודל להסביר לנו איך מכינים הרואין או שודדים בנק ונמלטים מהעונש עם אם מזיזים פלט של שכבה אחת במודל שפה. וגם ניתן למנוע ממודל ״לא מרוסן״ לתת תשובות לא פוגעניות ולפעמים להימנע מלענות על שאלות מסוכנות אם מזיזים את הפלטים של כל השכבות של מודל, כל אחת עם וקטור r_l כאשר l זו מספר השכבה. ", 2)

# 4) Second body paragraph
$d.Content.Find.Execute("אבל למה זה עובד בכלל? המחברים טוענים (באופן די מפתיע, אני חייב להגיד) שהביצועים הנפלאים של הטרנספורמרים נובעים בחלקם מיכולתה של פונקציית סופטמקס לכפות רגולריזציה מסוימת על נורמת פרובניוס של מטריצה המשקלים וגם של היעקוביאן שלה (ביחס לקלט של הסופטמקס) במהלך האימון הוא מסדר (sqrt(n כאשר n הינו מימד לקלט. ", $true, $false, $false, $false, $false, $true, 1, $false, "איך בעצם  מוצאים את הוקטורים האלה? עבור דאטהסט המכיל שאלות ותשובות רצויות מחשבים את ההפרש r הממוצע (על כל התשובות) בין האקטיבציות של כל שכבות המודל ועבור כל הטוקנים של חלון ההקשר. כלומר יש לנו מטריצה LxI של וקטורי ההפרש כאשר L זה מספר השכבות ו I זה מספר הטוקנים בחלון ההקשר.", 2)

# 5) Third body paragraph (contains a literal "r^T" -> escape the caret as "^^")
$d.Content.Find.Execute("נורמת פרובניוס או NF מוגדרת בתור שורש של סכום הריבועים של כל הערכים במטריצה והיא גם שווה לשורש של סכום הריבועים הערכים הסינגולריים (הכללה של ערכים עצמיים למטריצות לא ריבועיות). ד״א סופטמקס מחושב במנגנון ה-attention של מערך של וקטורים אז היעקוביאן תיאורטית הוא טנזור תלת מימדי (המאמר מפרט איך מחשבים את NF במקרה הזה).", $true, $false, $false, $false, $false, $true, 1, $false, "כדי לגרום למודל להיות ״פחות מרוסן״ אנו בוחרים שכבה שהוספתן של מורידה ממנו את בלמים בצורה המשמעותית ביותר (יש מדדים לא רעים לכך). כלומר משאירים I וקטורי הפרשים שחישבנו. כדי לגרום למודל להיות יותר מנומס צריך להחסיר את ״כיוון הגסות״ מכל השכבות של המודל בצורה שתעביר אותם ממרחב אורתוגונלי ל r (כל שכבה ולכל טוקן בחלון ההקשר). בפרט מכל אקטיבציה x בכל שכבה ובכל טוקן : r * r^^T *x קל לראות שהווקטור המתקבל כתוצאה מכך יהיה אורתוגונלי ל r. ", 2)

# 6) Fourth body paragraph
$d.Content.Find.Execute("אז בגדול המאמר מוכיח שני משפטים. בראשון מהם טוענים ש NF של מנגנון attention פולינומיאלי (כולל הלינארי) מתנהג לפי (O(n אם המטריצות שם, K ו-Q וגם ייצוגי הטוקנים מפולגים גאוסית כמובן). אז אם מנרמלים את ה-attention הפולינומיאלי עם (n^(-0.5 מקבלים את (sqrt(n שהיה לנו עבור מנגנון ה-attention הרגיל. בנוסף NF של היעקוביאן לפי Q, המנורמל לפי (n^(-0.5 (לא זה שמתנהג לפי (sqrt(n  ב-attention הרגיל) גם מתנהג לפני (sqrt(n.", $true, $false, $false, $false, $false, $true, 1, $false, "עושים זאת לווקטור האקטיבציה לפני residual connection בכל בלוק של טרנספורמר. כמובן (מכיוון שיש הרבה מכפלות של מטריצות)ניתן להזיז גם את המשקלים שלהם כדי לקבל את אותם האפקטים. מאמר די מגניב וקל להבנה.", 2)

# 7) Remove the now-obsolete closing paragraph entirely (the one right before the URL paragraph)
$target = "המחברים טוענים שזה מספיק כדי לטעון שניתן להחליף סופטמקס בפולינומים שיותר קלים מבחינה חישובית, מקבלים תוצאות מעודדות אבל אני עדיין לא השתכנעתי…"
$removed = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $p.Range.Delete()
        $removed = $true
        break
    }
}
if (-not $removed) {
    throw "Could not find the obsolete closing paragraph to delete"
}

# 8) URL paragraph
$d.Content.Find.Execute("https://arxiv.org/abs/2410.18613", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2406.11717", 2)
